$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 411.3846
$ws.Range("J19").Value = 453.7
$ws.Range("L19").Value = 453.7
$ws.Range("N19").Value = -803.7
$ws.Range("H63").Value = 59658.5
$ws.Range("J63").Value = 59071
$ws.Range("L63").Value = 59071
$ws.Range("N63").Value = -60319
$ws.Range("H66").Value = 59658.5
$ws.Range("J66").Value = 59071
$ws.Range("L66").Value = 177213
$ws.Range("N66").Value = -183453
$ws.Range("H75").Value = 60000
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61872
$ws.Range("H78").Value = 60000
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -189360
$ws.Range("H103").Value = 762.2727
$ws.Range("I103").Value = 263.5
$ws.Range("J103").Value = 873.1111
$ws.Range("K103").Value = 790.5
$ws.Range("L103").Value = 2619.3333
$ws.Range("M103").Value = -204.5
$ws.Range("N103").Value = -3791.3333
$ws.Range("H107").Value = 772.4762
$ws.Range("I107").Value = 772.4762
$ws.Range("K107").Value = 772.4762
$ws.Range("M107").Value = 1147.5238
$ws.Range("H111").Value = 13123.44
$ws.Range("I111").Value = 17222.357
$ws.Range("K111").Value = 51667.071
$ws.Range("M111").Value = -48600.071
$ws.Range("H112").Value = 386591.2
$ws.Range("J112").Value = 456780.5
$ws.Range("L112").Value = 1370341.5
$ws.Range("N112").Value = -1372557.5
$ws.Range("H138").Value = 163130.67
$ws.Range("J138").Value = 234749.48
$ws.Range("L138").Value = 704248.4400000001
$ws.Range("N138").Value = -714528.4400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10096.17
$ws.Range("I32").Value = 10096.17
$ws.Range("K32").Value = 10096.17
$ws.Range("M32").Value = -9809.17
$ws.Range("H45").Value = 3951.6086
$ws.Range("I45").Value = 3382.6667
$ws.Range("K45").Value = 3382.6667
$ws.Range("M45").Value = -3005.6667
$ws.Range("H63").Value = 1986.0769
$ws.Range("I63").Value = 1388.5714
$ws.Range("K63").Value = 1388.5714
$ws.Range("M63").Value = -702.5714
$ws.Range("H66").Value = 1986.0769
$ws.Range("I66").Value = 1388.5714
$ws.Range("K66").Value = 6942.857
$ws.Range("M66").Value = -3510.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2388937
$ws.Range("I86").Value = 5566947
$ws.Range("J86").Value = 5429.3125
$ws.Range("K86").Value = 5566947
$ws.Range("L86").Value = 5429.3125
$ws.Range("M86").Value = -5565824
$ws.Range("N86").Value = -7675.3125
$ws.Range("H89").Value = 2388937
$ws.Range("I89").Value = 5566947
$ws.Range("J89").Value = 5429.3125
$ws.Range("K89").Value = 27834735
$ws.Range("L89").Value = 27146.5625
$ws.Range("M89").Value = -27829119
$ws.Range("N89").Value = -38378.5625
$ws.Range("H105").Value = 6572.067
$ws.Range("I105").Value = 5758
$ws.Range("K105").Value = 5758
$ws.Range("M105").Value = -4011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3691.2354
$ws.Range("I31").Value = 2887.9167
$ws.Range("K31").Value = 2887.9167
$ws.Range("M31").Value = -2592.9167
$ws.Range("H34").Value = 3691.2354
$ws.Range("I34").Value = 2887.9167
$ws.Range("K34").Value = 2887.9167
$ws.Range("M34").Value = -2685.9167
$ws.Range("H97").Value = 84654
$ws.Range("J97").Value = 84654
$ws.Range("L97").Value = 84654
$ws.Range("N97").Value = -86636
$ws.Range("H122").Value = 4907.143
$ws.Range("J122").Value = 4751
$ws.Range("L122").Value = 14253
$ws.Range("N122").Value = -19153
$ws.Range("H132").Value = 1431341
$ws.Range("I132").Value = 1907433.8
$ws.Range("J132").Value = 3062.8572
$ws.Range("K132").Value = 5722301.4
$ws.Range("L132").Value = 9188.5716
$ws.Range("M132").Value = -5719771.4
$ws.Range("N132").Value = -14248.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6433493.5
$ws.Range("I4").Value = 6484812.5
$ws.Range("J4").Value = 6346250
$ws.Range("K4").Value = 19454437.5
$ws.Range("L4").Value = 19038750
$ws.Range("M4").Value = -19454325.5
$ws.Range("N4").Value = -19038974
$ws.Range("H80").Value = 3198.8
$ws.Range("J80").Value = 3249.75
$ws.Range("L80").Value = 9749.25
$ws.Range("N80").Value = -11621.25
$ws.Range("H83").Value = 3198.8
$ws.Range("J83").Value = 3249.75
$ws.Range("L83").Value = 29247.75
$ws.Range("N83").Value = -38607.75
$ws.Range("H113").Value = 1875.4286
$ws.Range("J113").Value = 1950.4615
$ws.Range("L113").Value = 5851.3845
$ws.Range("N113").Value = -10191.3845
$ws.Range("H122").Value = 1358
$ws.Range("I122").Value = 1499
$ws.Range("J122").Value = 1343.9
$ws.Range("K122").Value = 13491
$ws.Range("L122").Value = 12095.1
$ws.Range("M122").Value = -11041
$ws.Range("N122").Value = -16995.1
$ws.Range("H132").Value = 2630.55
$ws.Range("I132").Value = 1269.4
$ws.Range("J132").Value = 3447.24
$ws.Range("K132").Value = 11424.6
$ws.Range("L132").Value = 31025.16
$ws.Range("M132").Value = -8894.6
$ws.Range("N132").Value = -36085.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 98542.64
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 98542.64
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 98542.64
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -99512.64
$ws.Range("H102").Value = 13262.111
$ws.Range("I102").Value = 1779.3939
$ws.Range("J102").Value = 44839.582
$ws.Range("K102").Value = 1779.3939
$ws.Range("L102").Value = 44839.582
$ws.Range("M102").Value = -157.3939
$ws.Range("N102").Value = -48083.582
$ws.Range("H115").Value = 98542.64
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 98542.64
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 98542.64
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -100892.64
$ws.Range("H132").Value = 5837.304
$ws.Range("I132").Value = 5603.75
$ws.Range("J132").Value = 6371.143
$ws.Range("K132").Value = 16811.25
$ws.Range("L132").Value = 19113.429
$ws.Range("M132").Value = -14281.25
$ws.Range("N132").Value = -24173.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4164.875
$ws.Range("I40").Value = 3937.1724
$ws.Range("K40").Value = 3937.1724
$ws.Range("M40").Value = -3801.1724
$ws.Range("H46").Value = 11800.883
$ws.Range("I46").Value = 7366.6665
$ws.Range("J46").Value = 12751.071
$ws.Range("K46").Value = 7366.6665
$ws.Range("L46").Value = 12751.071
$ws.Range("M46").Value = -7178.6665
$ws.Range("N46").Value = -13127.071
$ws.Range("H93").Value = 1526.25
$ws.Range("I93").Value = 1393.125
$ws.Range("J93").Value = 1792.5
$ws.Range("K93").Value = 1393.125
$ws.Range("L93").Value = 1792.5
$ws.Range("M93").Value = -145.125
$ws.Range("N93").Value = -4288.5
$ws.Range("H100").Value = 3799.923
$ws.Range("I100").Value = 3299.875
$ws.Range("K100").Value = 3299.875
$ws.Range("M100").Value = -2758.875
$ws.Range("H122").Value = 3582.3215
$ws.Range("I122").Value = 3006.8
$ws.Range("K122").Value = 9020.400000000001
$ws.Range("M122").Value = -6570.400000000001
$ws.Range("H132").Value = 2449.3572
$ws.Range("I132").Value = 2247.2307
$ws.Range("K132").Value = 6741.6921
$ws.Range("M132").Value = -4211.6921
$ws.Range("H136").Value = 4365.057
$ws.Range("I136").Value = 3812.7827
$ws.Range("J136").Value = 5423.5835
$ws.Range("K136").Value = 11438.3481
$ws.Range("L136").Value = 16270.7505
$ws.Range("M136").Value = -8888.348100000001
$ws.Range("N136").Value = -21370.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2812.7778
$ws.Range("I96").Value = 1726.8334
$ws.Range("J96").Value = 3355.75
$ws.Range("K96").Value = 1726.8334
$ws.Range("L96").Value = 3355.75
$ws.Range("M96").Value = -353.8334
$ws.Range("N96").Value = -6101.75
$ws.Range("H100").Value = 1048.5555
$ws.Range("I100").Value = 1094.2858
$ws.Range("K100").Value = 2188.5716
$ws.Range("M100").Value = -1647.5716
$ws.Range("H132").Value = 2003.0385
$ws.Range("I132").Value = 2131.5908
$ws.Range("K132").Value = 6394.7724
$ws.Range("M132").Value = -3864.7724
